$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Part 1: the empty paragraph right after "... kontrol eder." gets the
# "$ dotnet new sln -n Day_01092024 Solition ile proje oluşturma" line,
# and the _GoBack bookmark moves here (it is a singleton bookmark, so
# re-adding it elsewhere removes the old one automatically).
# ------------------------------------------------------------------
$pCmd = $d.Paragraphs.Item(6)

# Borrow the existing "Lucida Console" 9pt run formatting (used by the
# other shell-command lines in this document) by copying a tiny bit of
# already-formatted text, then appending the rest of the line so it
# merges into one run with the same formatting.
$pSample = $d.Paragraphs.Item(2)
$sampleRange = $d.Range($pSample.Range.Start, $pSample.Range.Start + 2)
$pCmd.Range.FormattedText = $sampleRange.FormattedText
$pCmd.Range.InsertAfter("dotnet new sln -n Day_01092024 Solition ile proje oluşturma")

# Wrap just the text (not the paragraph mark) with the _GoBack bookmark.
$bmRange = $d.Range($pCmd.Range.Start, $pCmd.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# Part 2: the paragraph that used to hold the bookmark becomes a plain
# empty paragraph, and a new empty paragraph (same formatting) is
# inserted right after it.
# ------------------------------------------------------------------
$pOldBookmark = $d.Paragraphs.Item(7)
$pOldBookmark.Range.InsertParagraphAfter()

$pNew = $d.Paragraphs.Item(8)
$pNew.Range.InsertAfter("X")
$tmpRange = $d.Range($pNew.Range.Start, $pNew.Range.Start + 1)
$tmpRange.Delete()
